$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "First image " -> "First image  computation time = 0.5 sec"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("First image ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "First image  computation time = 0.5 sec", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: merge the "Second Image" paragraph with the following paragraph
# (tabs + spaces + picture), turning them into a single paragraph that keeps
# the numbered-list formatting of "Second Image" and whose text reads
# "Second Image      computation time = 0.5 sec" followed by the original
# tabs/spaces and the picture.
# ---------------------------------------------------------------------------
$secondImagePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "Second Image`r") {
        $secondImagePara = $candidate
        break
    }
}

$fullRange = $d.Range($secondImagePara.Range.Start, $secondImagePara.Range.End)

$xmlNamespaces = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' + `
    'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' + `
    'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' + `
    'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

$mergedParagraphXml = '<w:p ' + $xmlNamespaces + '>' + `
    '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:ind w:left="2160" w:hanging="360"/></w:pPr>' + `
    '<w:r><w:rPr><w:rtl w:val="0"/></w:rPr>' + `
    '<w:t xml:space="preserve">Second Image      computation time = 0.5 sec</w:t>' + `
    '<w:tab/><w:tab/><w:t xml:space="preserve">    </w:t></w:r>' + `
    '<w:r><w:drawing>' + `
    '<wp:anchor allowOverlap="1" behindDoc="0" distB="114300" distT="114300" distL="114300" distR="114300" hidden="0" layoutInCell="1" locked="0" relativeHeight="0" simplePos="0">' + `
    '<wp:simplePos x="0" y="0"/>' + `
    '<wp:positionH relativeFrom="column"><wp:posOffset>971550</wp:posOffset></wp:positionH>' + `
    '<wp:positionV relativeFrom="paragraph"><wp:posOffset>156163</wp:posOffset></wp:positionV>' + `
    '<wp:extent cx="5731200" cy="1422400"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:wrapNone/>' + `
    '<wp:docPr id="1" name="image3.png"/>' + `
    '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' + `
    '<pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image3.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr>' + `
    '<pic:blipFill><a:blip r:embed="rId8"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' + `
    '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5731200" cy="1422400"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr>' + `
    '</pic:pic></a:graphicData></a:graphic></wp:anchor>' + `
    '</w:drawing></w:r></w:p>'

$fullRange.InsertXML($mergedParagraphXml)

# Delete the now-duplicate paragraph that used to hold the tabs/spaces/picture
$leftoverPara = $d.Paragraphs.Item($secondImagePara.Index + 1)
$leftoverRange = $d.Range($leftoverPara.Range.Start, $leftoverPara.Range.End)
$leftoverRange.Delete()

# ---------------------------------------------------------------------------
# Change 3: after "Where k: corner sharpness = 0.04", insert two new
# paragraphs describing the pixel-thresholding step.
# ---------------------------------------------------------------------------
$cornerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "`tWhere k: corner sharpness = 0.04`r") {
        $cornerPara = $candidate
        break
    }
}

$cornerPara.Range.InsertParagraphAfter()
$firstNewPara = $d.Paragraphs.Item($cornerPara.Index + 1)
$firstNewTextRange = $d.Range($firstNewPara.Range.Start, $firstNewPara.Range.End - 1)
$firstNewTextRange.Text = "Get pixels above specific threshold and color it on the "

$firstNewParaFresh = $d.Paragraphs.Item($cornerPara.Index + 1)
$firstNewParaFresh.Range.InsertParagraphAfter()
$secondNewPara = $d.Paragraphs.Item($firstNewParaFresh.Index + 1)
$secondNewTextRange = $d.Range($secondNewPara.Range.Start, $secondNewPara.Range.End - 1)
$secondNewTextRange.Text = "    original image"
